$d = $word.ActiveDocument
$q = [char]34

# Table 1: ship particulars (2 columns, row label / value)
$t1 = $d.Tables.Item(1)

# Row 1: Ship name  "СИНЕГОРСК" -> "ВОЛГО-БАЛТ 136"  (value kept quoted)
$t1.Cell(1, 2).Range.Text = $q + "ВОЛГО-БАЛТ 136" + $q

# Row 2: Registry number  021026 -> 703999
$t1.Cell(2, 2).Range.Text = "703999"

# Row 3: Build date  22.07.2004 -> 22.03.1971
$t1.Cell(3, 2).Range.Text = "22.03.1971"

# Row 7: Gross tonnage  9611 -> 2457
$t1.Cell(7, 2).Range.Text = "2457"

# Row 8: Power, kW  5400 -> 1030
$t1.Cell(8, 2).Range.Text = "1030"

# Table 2: "1 Этап. ОФОРМЛЕНИЕ ДОКУМЕНТОВ" signatories table.
# Only the "механик" (row 3) and "эл. радио" (row 4) rows change name;
# the "корпус" (row 2) row keeps "Козлов С. В." unchanged.
$t2 = $d.Tables.Item(2)
$t2.Cell(3, 3).Range.Text = "Кудрявцев М. А."
$t2.Cell(4, 3).Range.Text = "Кудрявцев М. А."
